# Thread test & Annotation check
# Update the fully-qualified PageObject class name used in the "Data" sheet
# (MultipleThreadCacheTest -> file.pagefactory.excel.ExcelFPMultipleThreadCacheTest)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$newClassName = "file.pagefactory.excel.ExcelFPMultipleThreadCacheTest`$PageObjectExcelFirst"

$ws.Range("A2").Value = $newClassName
$ws.Range("A3").Value = $newClassName

# Widen column A so the longer class name keeps fitting (bestFit recalculation)
$ws.Columns.Item(1).ColumnWidth = 85.6

# Move the active selection on the Data sheet to B11
$ws.Range("B11").Select() | Out-Null
